# Auto-generated update of cryptos list data (Price / Volume(1h) columns)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '96.662.00'
$ws.Range('E2').Value = '  -1.57%  '
$ws.Range('D3').Value = '3.674.71'
$ws.Range('E3').Value = '  +1.50%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '239.78'
$ws.Range('E5').Value = '  -1.83%  '
$ws.Range('E6').Value = '  +7.50%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '655.55'
$ws.Range('E7').Value = '  -0.58%  '
$ws.Range('E8').Value = '  +0.26%  '
$ws.Range('E9').Value = '  +0.79%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.00'
$ws.Range('E10').Value = '  +0.07%  '
$ws.Range('D11').Value = '3.671.94'
$ws.Range('E11').Value = '  +1.37%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '45.50'
$ws.Range('E12').Value = '  +3.20%  '
$ws.Range('E13').Value = '  -0.15%  '
$ws.Range('E14').Value = '  +4.87%  '
$ws.Range('D15').Value = '4.359.11'
$ws.Range('E15').Value = '  +1.57%  '
$ws.Range('E16').Value = '  +2.79%  '
$ws.Range('D17').Value = '96.358.57'
$ws.Range('E17').Value = '  -1.71%  '
$ws.Range('D18').Value = '3.673.11'
$ws.Range('E18').Value = '  +1.51%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '18.78'
$ws.Range('E19').Value = '  +3.55%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.79'
$ws.Range('E20').Value = '  -1.18%  '
$ws.Range('E21').Value = '  -1.82%  '
$ws.Range('E22').Value = '  -1.56%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '529.84'
$ws.Range('E23').Value = '  +2.78%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.45'
$ws.Range('E24').Value = '  -1.66%  '
$ws.Range('E25').Value = '  +3.02%  '
$ws.Range('E26').Value = '  -2.20%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '101.68'
$ws.Range('E27').Value = '  +1.76%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '13.19'
$ws.Range('E28').Value = '  +1.34%  '
$ws.Range('D29').Value = '3.870.04'
$ws.Range('E29').Value = '  +1.47%  '
$ws.Range('E30').Value = '  +8.70%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '12.45'
$ws.Range('E31').Value = '  +5.30%  '
$ws.Range('E32').Value = '  -0.52%  '
$ws.Range('E33').Value = '  +0.05%  '
$ws.Range('E34').Value = '  +15.63%  '
$ws.Range('E35').Value = '  -0.78%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '673.86'
$ws.Range('E36').Value = '  +10.04%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.00'
$ws.Range('E37').Value = '  +0.65%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '32.41'
$ws.Range('E38').Value = '  +1.52%  '
$ws.Range('E39').Value = '  +3.12%  '
$ws.Range('E40').Value = '  -0.71%  '
$ws.Range('E41').Value = '  +4.61%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.99'
$ws.Range('E42').Value = '  -0.36%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.959'
$ws.Range('E43').Value = '  +3.03%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '6.52'
$ws.Range('E44').Value = '  +8.64%  '
$ws.Range('E45').Value = '  +17.37%  '
$ws.Range('E46').Value = '  +0.06%  '
$ws.Range('E47').Value = '  +5.57%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.444'
$ws.Range('E48').Value = '  +11.46%  '
$ws.Range('E49').Value = '  +0.61%  '
$ws.Range('E50').Value = '  +4.89%  '
$ws.Range('E51').Value = '  -0.02%  '

# Restore default (Normal) style on cells where we forced a Text number format,
# so the resulting style matches the workbook's original default formatting.
$ws.Range('D5').Style = 'Normal'
$ws.Range('D7').Style = 'Normal'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D31').Style = 'Normal'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D48').Style = 'Normal'
